$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '41.733.43'
Set-TextValue 'E2' '  -1.87%  '
Set-TextValue 'D3' '2.218.06'
Set-TextValue 'E3' '  -1.59%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '241.58'
Set-TextValue 'E5' '  -2.08%  '
Set-TextValue 'E6' '  -1.03%  '
Set-TextValue 'D7' '72.64'
Set-TextValue 'E7' '  -5.62%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.596'
Set-TextValue 'E9' '  -4.38%  '
Set-TextValue 'D10' '42.16'
Set-TextValue 'E10' '  -6.87%  '
Set-TextValue 'D11' '0.0951'
Set-TextValue 'E11' '  -0.12%  '
Set-TextValue 'E12' '  -4.56%  '
Set-TextValue 'D13' '0.104'
Set-TextValue 'E13' '  +0.51%  '
Set-TextValue 'D14' '2.551.67'
Set-TextValue 'E14' '  -1.62%  '
Set-TextValue 'D15' '14.27'
Set-TextValue 'E15' '  -2.64%  '
Set-TextValue 'D16' '0.836'
Set-TextValue 'E16' '  -2.98%  '
Set-TextValue 'D17' '2.219.42'
Set-TextValue 'E17' '  -3.15%  '
Set-TextValue 'D18' '41.659.70'
Set-TextValue 'D19' '0.0000106'
Set-TextValue 'E19' '  +3.81%  '
Set-TextValue 'D20' '72.51'
Set-TextValue 'E20' '  +0.42%  '
Set-TextValue 'D21' '6.16'
Set-TextValue 'E21' '  -0.80%  '
Set-TextValue 'D22' '10.98'
Set-TextValue 'E22' '  +20.60%  '
Set-TextValue 'D23' '229.47'
Set-TextValue 'E23' '  -1.06%  '
Set-TextValue 'D24' '2.06'
Set-TextValue 'E24' '  -8.98%  '
Set-TextValue 'E25' '  +0.14%  '
Set-TextValue 'D26' '11.41'
Set-TextValue 'E26' '  -1.43%  '
Set-TextValue 'D27' '3.62'
Set-TextValue 'E27' '  -0.02%  '
Set-TextValue 'E28' '  -1.74%  '
Set-TextValue 'E29' '  -0.66%  '
Set-TextValue 'D30' '167.20'
Set-TextValue 'E30' '  -0.43%  '
Set-TextValue 'D31' '20.46'
Set-TextValue 'E31' '  -1.24%  '
Set-TextValue 'E32' '  -3.84%  '
Set-TextValue 'D33' '5.51'
Set-TextValue 'E33' '  +3.73%  '
Set-TextValue 'D34' '30.11'
Set-TextValue 'E34' '  -3.83%  '
Set-TextValue 'E35' '  -0.56%  '
Set-TextValue 'E36' '  -10.27%  '
Set-TextValue 'E37' '  -6.21%  '
Set-TextValue 'E38' '  -3.99%  '
Set-TextValue 'D39' '13.33'
Set-TextValue 'E39' '  -6.25%  '
Set-TextValue 'E40' '  -2.92%  '
Set-TextValue 'D41' '64.35'
Set-TextValue 'E41' '  +0.55%  '
Set-TextValue 'D42' '5.61'
Set-TextValue 'E42' '  -3.73%  '
Set-TextValue 'E43' '  -2.63%  '
Set-TextValue 'E44' '  -1.20%  '
Set-TextValue 'D45' '103.09'
Set-TextValue 'E45' '  -4.80%  '
Set-TextValue 'D46' '0.100'
Set-TextValue 'E46' '  -2.97%  '
Set-TextValue 'D47' '2.33'
Set-TextValue 'E47' '  -2.61%  '
Set-TextValue 'B48' 'ARBITRUM'
Set-TextValue 'C48' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D48' '1.11'
Set-TextValue 'E48' '  -2.63%  '
Set-TextValue 'B49' 'TrustWalletToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D49' '1.16'
Set-TextValue 'E49' '  -2.54%  '
Set-TextValue 'E50' '  -0.74%  '
Set-TextValue 'D51' '2.423.64'
Set-TextValue 'E51' '  -1.80%  '
